$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A3").Value = "Pstream"
$ws.Range("B3").Value = "https://pstream.org/"
$ws.Range("C3").Value = "PStream.jpeg"

$excel.ActiveWindow.ScrollRow = 13
[void]$ws.Range("B18").Select()
